$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reuse the same text/style pattern as existing rows (col A uses the
# "vertical center" style already applied to A2:A20).
$ws.Range("A21").Value = "Deals_Chat_ShipperUser_TC001"
$ws.Range("B21").Value = "John Tucker"
$ws.Range("C21").Value = "ONE"
$ws.Range("D21").Value = "Deal shared successfully"

$ws.Range("A22").Value = "Deals_Chat_ShipperAdmin_TC002"
$ws.Range("B22").Value = "John Tucker"
$ws.Range("C22").Value = "ONE"
$ws.Range("D22").Value = "Deal shared successfully"

$ws.Range("A23").Value = "Deals_Chat_CarrierUser_TC003"
$ws.Range("B23").Value = "Stan Koster Andersons"
$ws.Range("C23").Value = "ONE"
$ws.Range("D23").Value = "Deal shared successfully"

# Match style used by column A on existing data rows (vertical-center alignment)
$ws.Range("A21:A23").VerticalAlignment = -4108

# Update view: scroll so row 4 is at top, and select B23 as the active cell
$ws.Range("B23").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 4 | Out-Null
